$wb = $excel.ActiveWorkbook

# --- Rename sheets ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "Test Datatype Array"
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "Test Array"
$ws3 = $wb.Worksheets.Item(3)

# --- Sheet1 selection update ---
$ws1.Range("C7:E7").Select()

# --- Sheet2 content (new "Test Array" sheet) ---
$ws2.Range("C5:E5").Merge()
$ws2.Range("C5").Value = "Spreadsheet DoubleValue[] testDoubleValue()`n"
$ws2.Range("I5").Value = "Method DoubleValue[] getDVs()"

$ws2.Range("D6:E6").Merge()
$ws2.Range("D6").Value = "Value"
$ws2.Range("I6").Value = "return new DoubleValue[]{new DoubleValue(1.23), new DoubleValue(5.24)};"

$ws2.Range("C7").Value = "AnyValue"
$ws2.Range("D7").Value = "= 78"

$ws2.Range("C8").Value = ""
$ws2.Range("D8").Value = "{getDVs()}"

$ws2.Range("C9").Value = "RETURN"
$ws2.Range("D9:E9").Merge()
$ws2.Range("D9").Value = "'=`$DVs"

# Row heights
$ws2.Rows.Item(5).RowHeight = 25.5

# Column widths (best effort given engine's column-width rounding model)
$ws2.Columns.Item(3).ColumnWidth = 16
$ws2.Columns.Item(5).ColumnWidth = 33.16666666666667
$ws2.Columns.Item(9).ColumnWidth = 69.33333333333333

# Page setup
$ws2.PageSetup.PaperSize = 9
$ws2.PageSetup.Orientation = 1

# Selection on sheet2
$ws2.Range("C7").Select()

# --- Activate sheet2 (Test Array) as the active tab ---
$ws2.Activate()
